# Apply updated "想去人数" (interest count) / "最低票价" (min ticket price)
# figures to the 展览 (Exhibition) and 全部类型 (All Types) sheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F4").Value = 22
    $ws.Range("G4").Value = 55

    $ws.Range("F6").Value = 1179
    $ws.Range("G6").Value = 70

    $ws.Range("G7").Value = 70

    $ws.Range("G8").Value = 60

    if ($sheetName -eq "展览") {
        $ws.Range("F14").Value = 96
        $ws.Range("F17").Value = 289
        $ws.Range("F18").Value = 396
        $ws.Range("F19").Value = 488
        $ws.Range("F21").Value = 5953
        $ws.Range("F22").Value = 5290
    }
    else {
        $ws.Range("F16").Value = 96
        $ws.Range("F19").Value = 289
        $ws.Range("F20").Value = 396
        $ws.Range("F21").Value = 488
        $ws.Range("F23").Value = 5953
        $ws.Range("F25").Value = 5290
    }
}
